$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the insuredId value in the rule table (row 10, column B)
$ws.Range("B10").Value = "A223456789"

# Reflect the new active cell selection left by the author
$ws.Range("B10").Select()
